# Update "想去人数" (F column) values for several events on the
# "展览" and "全部类型" sheets, as output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value = 2361
$wsExhibit.Range("F13").Value = 1431
$wsExhibit.Range("F14").Value = 506
$wsExhibit.Range("F16").Value = 312
$wsExhibit.Range("F26").Value = 1475
$wsExhibit.Range("F28").Value = 371
$wsExhibit.Range("F29").Value = 228

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 2361
$wsAll.Range("F14").Value = 1431
$wsAll.Range("F15").Value = 506
$wsAll.Range("F17").Value = 312
$wsAll.Range("F27").Value = 1475
$wsAll.Range("F29").Value = 371
$wsAll.Range("F30").Value = 228
